$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates: B2 cleared, C2/D2/E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.97080645429713486
$ws.Range("D2").Value = 0.13896637662626507
$ws.Range("E2").Value = 3.2553633475097685

# Row 3 updates
$ws.Range("B3").Value = 0.22509082614623285
$ws.Range("C3").Value = 1.592273790255478
$ws.Range("D3").Value = 0.24152064206455273
$ws.Range("E3").Value = 1.3791337152805876

# Update selection to match new selection range
$ws.Range("B1:E3").Select()
